$wb = $excel.ActiveWorkbook

# Map of sheet index (1-based, tab order) -> list of (cellRef, expectedOldValue, newValue)
# Sheet1 = 展览 (Exhibition), Sheet2 = 演出 (Performance),
# Sheet3 = 本地生活 (Local Life), Sheet4 = 全部类型 (All Types)
$changesBySheet = @{}
$changesBySheet[1] = @(
    @{ Cell = "F3"; OldValue = 3830; NewValue = 3834 }
    @{ Cell = "F5"; OldValue = 1368; NewValue = 1369 }
    @{ Cell = "F6"; OldValue = 3828; NewValue = 3831 }
    @{ Cell = "F7"; OldValue = 389; NewValue = 390 }
    @{ Cell = "F8"; OldValue = 196; NewValue = 198 }
    @{ Cell = "F10"; OldValue = 8610; NewValue = 8639 }
    @{ Cell = "F11"; OldValue = 485; NewValue = 489 }
    @{ Cell = "F13"; OldValue = 133; NewValue = 134 }
    @{ Cell = "F14"; OldValue = 116; NewValue = 123 }
    @{ Cell = "F15"; OldValue = 291; NewValue = 293 }
    @{ Cell = "F16"; OldValue = 334; NewValue = 335 }
    @{ Cell = "F17"; OldValue = 91; NewValue = 92 }
    @{ Cell = "F18"; OldValue = 355; NewValue = 359 }
    @{ Cell = "F19"; OldValue = 10894; NewValue = 10946 }
    @{ Cell = "F21"; OldValue = 142; NewValue = 143 }
    @{ Cell = "F22"; OldValue = 389; NewValue = 390 }
    @{ Cell = "F23"; OldValue = 185; NewValue = 187 }
    @{ Cell = "F24"; OldValue = 10; NewValue = 11 }
    @{ Cell = "F26"; OldValue = 91; NewValue = 92 }
    @{ Cell = "F27"; OldValue = 136; NewValue = 137 }
    @{ Cell = "F28"; OldValue = 2679; NewValue = 2680 }
    @{ Cell = "F29"; OldValue = 2077; NewValue = 2080 }
    @{ Cell = "F30"; OldValue = 37; NewValue = 39 }
    @{ Cell = "F32"; OldValue = 2119; NewValue = 2120 }
    @{ Cell = "F33"; OldValue = 894; NewValue = 896 }
    @{ Cell = "F34"; OldValue = 4080; NewValue = 4085 }
    @{ Cell = "F35"; OldValue = 2564; NewValue = 2567 }
    @{ Cell = "F36"; OldValue = 279; NewValue = 281 }
    @{ Cell = "F37"; OldValue = 2582; NewValue = 2585 }
    @{ Cell = "F38"; OldValue = 3026; NewValue = 3027 }
    @{ Cell = "F39"; OldValue = 1246; NewValue = 1247 }
    @{ Cell = "F40"; OldValue = 170; NewValue = 172 }
    @{ Cell = "F42"; OldValue = 341; NewValue = 344 }
    @{ Cell = "F43"; OldValue = 319; NewValue = 325 }
    @{ Cell = "F44"; OldValue = 43; NewValue = 44 }
    @{ Cell = "F45"; OldValue = 107; NewValue = 110 }
    @{ Cell = "F46"; OldValue = 125; NewValue = 126 }
    @{ Cell = "F48"; OldValue = 96; NewValue = 98 }
    @{ Cell = "F49"; OldValue = 83; NewValue = 84 }
)

$changesBySheet[2] = @(
    @{ Cell = "F2"; OldValue = 15; NewValue = 16 }
    @{ Cell = "F7"; OldValue = 46; NewValue = 47 }
    @{ Cell = "F20"; OldValue = 1; NewValue = 2 }
    @{ Cell = "F22"; OldValue = 42; NewValue = 50 }
)

$changesBySheet[3] = @(
    @{ Cell = "F3"; OldValue = 33; NewValue = 34 }
)

$changesBySheet[4] = @(
    @{ Cell = "F3"; OldValue = 3830; NewValue = 3834 }
    @{ Cell = "F6"; OldValue = 1368; NewValue = 1369 }
    @{ Cell = "F7"; OldValue = 3828; NewValue = 3831 }
    @{ Cell = "F8"; OldValue = 389; NewValue = 390 }
    @{ Cell = "F10"; OldValue = 196; NewValue = 198 }
    @{ Cell = "F11"; OldValue = 8610; NewValue = 8639 }
    @{ Cell = "F12"; OldValue = 485; NewValue = 489 }
    @{ Cell = "F13"; OldValue = 116; NewValue = 123 }
    @{ Cell = "F14"; OldValue = 291; NewValue = 293 }
    @{ Cell = "F15"; OldValue = 334; NewValue = 335 }
    @{ Cell = "F16"; OldValue = 91; NewValue = 92 }
    @{ Cell = "F17"; OldValue = 355; NewValue = 359 }
    @{ Cell = "F18"; OldValue = 10894; NewValue = 10947 }
    @{ Cell = "F20"; OldValue = 33; NewValue = 34 }
    @{ Cell = "F21"; OldValue = 142; NewValue = 143 }
    @{ Cell = "F22"; OldValue = 389; NewValue = 390 }
    @{ Cell = "F23"; OldValue = 185; NewValue = 187 }
    @{ Cell = "F24"; OldValue = 10; NewValue = 11 }
    @{ Cell = "F27"; OldValue = 91; NewValue = 92 }
    @{ Cell = "F28"; OldValue = 136; NewValue = 137 }
    @{ Cell = "F29"; OldValue = 2679; NewValue = 2680 }
    @{ Cell = "F30"; OldValue = 2077; NewValue = 2080 }
    @{ Cell = "F31"; OldValue = 2119; NewValue = 2120 }
    @{ Cell = "F32"; OldValue = 894; NewValue = 896 }
    @{ Cell = "F33"; OldValue = 1; NewValue = 2 }
    @{ Cell = "F34"; OldValue = 4080; NewValue = 4085 }
    @{ Cell = "F35"; OldValue = 2564; NewValue = 2567 }
    @{ Cell = "F36"; OldValue = 279; NewValue = 281 }
    @{ Cell = "F37"; OldValue = 2582; NewValue = 2585 }
    @{ Cell = "F38"; OldValue = 3026; NewValue = 3027 }
    @{ Cell = "F39"; OldValue = 42; NewValue = 50 }
    @{ Cell = "F40"; OldValue = 1246; NewValue = 1247 }
    @{ Cell = "F41"; OldValue = 170; NewValue = 172 }
    @{ Cell = "F43"; OldValue = 341; NewValue = 344 }
    @{ Cell = "F44"; OldValue = 319; NewValue = 325 }
    @{ Cell = "F45"; OldValue = 107; NewValue = 110 }
    @{ Cell = "F46"; OldValue = 125; NewValue = 126 }
    @{ Cell = "F48"; OldValue = 96; NewValue = 98 }
    @{ Cell = "F49"; OldValue = 83; NewValue = 84 }
)

$totalApplied = 0
$mismatches = 0
foreach ($sheetIdx in $changesBySheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    foreach ($change in $changesBySheet[$sheetIdx]) {
        $cell = $ws.Range($change.Cell)
        $current = $cell.Value2
        if ($current -ne $change.OldValue) {
            $mismatches++
            Write-Output ("MISMATCH sheet=" + $sheetIdx + " cell=" + $change.Cell + " expectedOld=" + $change.OldValue + " actualOld=" + $current)
        }
        $cell.Value = $change.NewValue
        $totalApplied++
    }
}
Write-Output ("Applied " + $totalApplied + " changes with " + $mismatches + " mismatches")
